# ---------------------------------------------------------------------------
# B1--and-B2-PowerPoint.pptx edit
#
# 1. Slide 5 contains a 3x6 table (the "Type of document" table). Its table
#    style is switched from the default Table_0 style
#    ({58E01270-4AA8-47C4-940C-32E25577ADF0}) to the built-in "Medium Style 2
#    - Accent 1" gallery style ({1A427C4B-1E59-4C71-B9DB-1A960CFC5637}).
#
# 2. The presentation's Design/Theme is switched from "Integral" (Red
#    Violet colour scheme) to the plain "Office Theme" (standard Office blue
#    colour scheme). We reproduce this by writing the 12 standard Office
#    theme colours into the presentation's live theme colour scheme.
# ---------------------------------------------------------------------------

$p = $ppt.ActivePresentation

# --- 1. Table style on slide 5 ------------------------------------------------
$tableSlide = $p.Slides.Item(5)
$tableShape = $tableSlide.Shapes.Item(2)
$tableShape.Table.ApplyStyle("{1A427C4B-1E59-4C71-B9DB-1A960CFC5637}")

# --- 2. Swap the active theme's colour scheme: Integral -> Office Theme ------
# (dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink - in that canonical order)
$colors = $p.Slides.Item(1).ThemeColorScheme
$colors.Item(1).RGB  = 0          # dk1      000000
$colors.Item(2).RGB  = 16777215   # lt1      FFFFFF
$colors.Item(3).RGB  = 6968388    # dk2      44546A
$colors.Item(4).RGB  = 15132391   # lt2      E7E6E6
$colors.Item(5).RGB  = 13998939   # accent1  5B9BD5
$colors.Item(6).RGB  = 3243501    # accent2  ED7D31
$colors.Item(7).RGB  = 10855845   # accent3  A5A5A5
$colors.Item(8).RGB  = 49407      # accent4  FFC000
$colors.Item(9).RGB  = 12874308   # accent5  4472C4
$colors.Item(10).RGB = 4697456    # accent6  70AD47
$colors.Item(11).RGB = 12673797   # hlink    0563C1
$colors.Item(12).RGB = 7491477    # folHlink 954F72
